$wb = $excel.ActiveWorkbook

# ===== Sheet "LP1912": refresh scrape timestamp, row count, and schedule rows =====
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 04:17:50"
$ws1.Cells.Item(3,1).Value = "Total filas: 10"

$ws1.Cells.Item(6,1).Value = "04:17:50"
$ws1.Cells.Item(6,2).Value = "04:45"
$ws1.Cells.Item(6,3).Value = "215A_EL PATO"
$ws1.Cells.Item(6,4).Value = 28
$ws1.Cells.Item(6,5).Value = "LP1912"

$ws1.Cells.Item(7,1).Value = "04:17:50"
$ws1.Cells.Item(7,2).Value = "04:53"
$ws1.Cells.Item(7,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(7,4).Value = 36
$ws1.Cells.Item(7,5).Value = "LP1912"

$ws1.Cells.Item(8,1).Value = "04:17:50"
$ws1.Cells.Item(8,2).Value = "05:16"
$ws1.Cells.Item(8,3).Value = "17_ROMERO"
$ws1.Cells.Item(8,4).Value = 59
$ws1.Cells.Item(8,5).Value = "LP1912"

$ws1.Cells.Item(9,1).Value = "04:17:50"
$ws1.Cells.Item(9,2).Value = "05:22"
$ws1.Cells.Item(9,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(9,4).Value = 65
$ws1.Cells.Item(9,5).Value = "LP1912"

$ws1.Cells.Item(10,1).Value = "04:17:50"
$ws1.Cells.Item(10,2).Value = "05:34"
$ws1.Cells.Item(10,3).Value = "215B_EL PATO"
$ws1.Cells.Item(10,4).Value = 77
$ws1.Cells.Item(10,5).Value = "LP1912"

$ws1.Cells.Item(11,1).Value = "04:17:50"
$ws1.Cells.Item(11,2).Value = "05:46"
$ws1.Cells.Item(11,3).Value = "15_ABASTO"
$ws1.Cells.Item(11,4).Value = 89
$ws1.Cells.Item(11,5).Value = "LP1912"

$ws1.Cells.Item(12,1).Value = "04:17:50"
$ws1.Cells.Item(12,2).Value = "05:54"
$ws1.Cells.Item(12,3).Value = "10_OLMOS"
$ws1.Cells.Item(12,4).Value = 97
$ws1.Cells.Item(12,5).Value = "LP1912"

$ws1.Cells.Item(13,1).Value = "04:17:50"
$ws1.Cells.Item(13,2).Value = "06:04"
$ws1.Cells.Item(13,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(13,4).Value = 107
$ws1.Cells.Item(13,5).Value = "LP1912"

$ws1.Cells.Item(14,1).Value = "04:17:50"
$ws1.Cells.Item(14,2).Value = "06:11"
$ws1.Cells.Item(14,3).Value = "215A_EL PATO"
$ws1.Cells.Item(14,4).Value = 114
$ws1.Cells.Item(14,5).Value = "LP1912"

$ws1.Cells.Item(15,1).Value = "04:17:50"
$ws1.Cells.Item(15,2).Value = "06:14"
$ws1.Cells.Item(15,3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(15,4).Value = 117
$ws1.Cells.Item(15,5).Value = "LP1912"

# ===== Sheet "LP1912-215": rebuild with the refreshed 215-filtered subset =====
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Clear()
$ws2.Cells.Item(1,1).Value = "LÍNEA 141 - LP1912-215 - 31/01/2026"
$ws2.Cells.Item(2,1).Value = "Última actualización: 04:17:50"
$ws2.Cells.Item(3,1).Value = "Total filas: 3"

$ws2.Cells.Item(5,1).Value = "Hora_Scrap"
$ws2.Cells.Item(5,2).Value = "Hora_Llegada"
$ws2.Cells.Item(5,3).Value = "Linea"
$ws2.Cells.Item(5,4).Value = "Minutos"
$ws2.Cells.Item(5,5).Value = "Parada"

$ws2.Cells.Item(6,1).Value = "04:17:50"
$ws2.Cells.Item(6,2).Value = "04:45"
$ws2.Cells.Item(6,3).Value = "215A_EL PATO"
$ws2.Cells.Item(6,4).Value = 28
$ws2.Cells.Item(6,5).Value = "LP1912"

$ws2.Cells.Item(7,1).Value = "04:17:50"
$ws2.Cells.Item(7,2).Value = "05:34"
$ws2.Cells.Item(7,3).Value = "215B_EL PATO"
$ws2.Cells.Item(7,4).Value = 77
$ws2.Cells.Item(7,5).Value = "LP1912"

$ws2.Cells.Item(8,1).Value = "04:17:50"
$ws2.Cells.Item(8,2).Value = "06:11"
$ws2.Cells.Item(8,3).Value = "215A_EL PATO"
$ws2.Cells.Item(8,4).Value = 114
$ws2.Cells.Item(8,5).Value = "LP1912"

# Re-apply the bold title/header style (same as sheet "LP1912") to rows 1-3
$ws1.Range("A1:A3").Copy()
$ws2.Range("A1:A3").PasteSpecial(-4122)

# ===== Sheet "6203-6173": refresh title/timestamp and rebuild header + schedule rows =====
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(1,1).Value = "LÍNEA 141 - 6203-6173 - 31/01/2026"
$ws3.Cells.Item(2,1).Value = "Última actualización: 04:17:50"
$ws3.Cells.Item(3,1).Value = "Total filas: 2"

$ws3.Cells.Item(5,1).Value = "Hora_Scrap"
$ws3.Cells.Item(5,2).Value = "Hora_Llegada"
$ws3.Cells.Item(5,3).Value = "Linea"
$ws3.Cells.Item(5,4).Value = "Minutos"
$ws3.Cells.Item(5,5).Value = "Parada"

$ws3.Cells.Item(6,1).Value = "04:17:50"
$ws3.Cells.Item(6,2).Value = "05:43"
$ws3.Cells.Item(6,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6,4).Value = 86
$ws3.Cells.Item(6,5).Value = "L6173"

$ws3.Cells.Item(7,1).Value = "04:17:50"
$ws3.Cells.Item(7,2).Value = "06:08"
$ws3.Cells.Item(7,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7,4).Value = 111
$ws3.Cells.Item(7,5).Value = "L6173"

